$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 3761
$ws.Range("K3").Value = 3794
$ws.Range("H4").Value = 1734
$ws.Range("J4").Value = 1820
$ws.Range("K4").Value = 774
$ws.Range("K5").Value = 262
$ws.Range("K6").Value = 4325
$ws.Range("H7").Value = 26047
$ws.Range("J7").Value = 29289
$ws.Range("K7").Value = 12916

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K3").Value = 38
$ws.Range("K6").Value = 89
$ws.Range("K7").Value = 176

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K3").Value = 262
$ws.Range("K6").Value = 286
$ws.Range("K7").Value = 869

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 201
$ws.Range("K6").Value = 152
$ws.Range("K7").Value = 531

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 106
$ws.Range("K5").Value = 26
$ws.Range("K6").Value = 98
$ws.Range("K7").Value = 372
$ws.Range("K8").Value = 869
$ws.Range("K15").Value = 129
$ws.Range("K20").Value = 292
$ws.Range("K21").Value = 33
$ws.Range("K27").Value = 128
$ws.Range("K29").Value = 683
$ws.Range("K33").Value = 531
$ws.Range("K40").Value = 32
$ws.Range("K42").Value = 455
$ws.Range("K44").Value = 120
$ws.Range("K46").Value = 30
$ws.Range("K47").Value = 74
$ws.Range("K49").Value = 73
$ws.Range("K50").Value = 72
$ws.Range("K51").Value = 151
$ws.Range("K52").Value = 355
$ws.Range("K53").Value = 176
$ws.Range("K54").Value = 249
$ws.Range("K59").Value = 24
$ws.Range("K60").Value = 83
$ws.Range("H63").Value = 286
$ws.Range("J63").Value = 105
$ws.Range("K63").Value = 45
$ws.Range("K64").Value = 78
$ws.Range("K67").Value = 511
$ws.Range("K73").Value = 116
$ws.Range("K76").Value = 188
$ws.Range("K78").Value = 160
$ws.Range("K79").Value = 328
$ws.Range("K84").Value = 94
$ws.Range("K85").Value = 587
$ws.Range("K88").Value = 147
$ws.Range("K89").Value = 179
$ws.Range("K90").Value = 118
$ws.Range("K91").Value = 139
$ws.Range("K96").Value = 153
$ws.Range("K98").Value = 66
$ws.Range("H101").Value = 26047
$ws.Range("J101").Value = 29289
$ws.Range("K101").Value = 12916

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 174
$ws.Range("K6").Value = 148
$ws.Range("K7").Value = 511

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K3").Value = 38
$ws.Range("K7").Value = 94

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K6").Value = 42
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K2").Value = 43
$ws.Range("K6").Value = 118
$ws.Range("K7").Value = 249

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K3").Value = 239
$ws.Range("K6").Value = 196
$ws.Range("K7").Value = 683

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K3").Value = 35
$ws.Range("K7").Value = 120

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K3").Value = 36
$ws.Range("K7").Value = 188

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K2").Value = 40
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 121
$ws.Range("K3").Value = 149
$ws.Range("K7").Value = 455

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value = 46
$ws.Range("K4").Value = 15
$ws.Range("K5").Value = 5
$ws.Range("K6").Value = 60
$ws.Range("K7").Value = 160

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K3").Value = 28
$ws.Range("K7").Value = 153

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K6").Value = 32
$ws.Range("K7").Value = 139

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 114
$ws.Range("K7").Value = 328

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K2").Value = 15
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 88
$ws.Range("K7").Value = 292

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 139
$ws.Range("K3").Value = 114
$ws.Range("K5").Value = 15
$ws.Range("K7").Value = 372

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 74

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 47
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 66

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 72

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 116

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 24

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K2").Value = 34
$ws.Range("K7").Value = 147

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K2").Value = 44
$ws.Range("K7").Value = 179

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 26

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K6").Value = 50
$ws.Range("K7").Value = 128

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 44
$ws.Range("K5").Value = 5
$ws.Range("K7").Value = 118

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K3").Value = 43
$ws.Range("K7").Value = 151

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K2").Value = 27
$ws.Range("K3").Value = 27
$ws.Range("K7").Value = 83

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K4").Value = 33
$ws.Range("K6").Value = 134
$ws.Range("K7").Value = 587

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K3").Value = 93
$ws.Range("K4").Value = 18
$ws.Range("K6").Value = 140
$ws.Range("K7").Value = 355
